# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") for rows 2-38, recalculated from the
# underlying strike data (replacing the old "Strike#" derived values).
$newK = @{
    2  = 2
    3  = 2
    4  = 2
    5  = 2
    6  = 2
    7  = 0
    8  = 2
    9  = 2
    10 = 4
    11 = 3
    12 = 1
    13 = 3
    14 = 1
    15 = 0
    16 = 2
    17 = 0
    18 = 2
    19 = 2
    20 = 2
    21 = 3
    22 = 0
    23 = 1
    24 = 2
    25 = 2
    26 = 0
    27 = 0
    28 = 0
    29 = 1
    30 = 0
    31 = 0
    32 = 1
    33 = 1
    34 = 1
    35 = 0
    36 = 1
    37 = 1
    38 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
